$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "all" -> new data row 38 (was the "footnote" row), footnote pushes to
# row 39.
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")
$wsAll.Rows.Item(38).Insert()

$wsAll.Cells.Item(38, 1).Value = 43966
$wsAll.Cells.Item(38, 1).NumberFormat = "m/d/yyyy"
$wsAll.Cells.Item(38, 2).Value = 281
$wsAll.Cells.Item(38, 3).Value = 277
$wsAll.Cells.Item(38, 4).Value = 64
$wsAll.Cells.Item(38, 5).Value = 54
$wsAll.Cells.Item(38, 6).Value = 10
$wsAll.Cells.Item(38, 7).Value = 11
$wsAll.Cells.Item(38, 8).Value = 202

$wsAll.Cells.Item(38, 1).Style = $wsAll.Cells.Item(37, 1).Style
$wsAll.Cells.Item(38, 2).Style = $wsAll.Cells.Item(37, 2).Style
$wsAll.Cells.Item(38, 3).Style = $wsAll.Cells.Item(37, 3).Style
$wsAll.Cells.Item(38, 4).Style = $wsAll.Cells.Item(37, 4).Style
$wsAll.Cells.Item(38, 5).Style = $wsAll.Cells.Item(37, 5).Style
$wsAll.Cells.Item(38, 6).Style = $wsAll.Cells.Item(37, 6).Style
$wsAll.Cells.Item(38, 7).Style = $wsAll.Cells.Item(37, 7).Style
$wsAll.Cells.Item(38, 8).Style = $wsAll.Cells.Item(37, 8).Style

$wsAll.Application.ActiveWindow.FreezePanes = $false
$wsAll.Range("E32").Select()
$wsAll.Application.ActiveWindow.FreezePanes = $true
$wsAll.Range("G42").Select()

# ---------------------------------------------------------------------------
# Sheet "kobe" -> new data row 93 (was the "footnote" row), footnote pushes
# to row 94. This sheet stops being the "tabSelected" sheet.
# ---------------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")
$wsKobe.Rows.Item(93).Insert()

$wsKobe.Cells.Item(93, 1).Value = 43966
$wsKobe.Cells.Item(93, 1).NumberFormat = "m/d/yyyy"
$wsKobe.Cells.Item(93, 2).Value = 0
$wsKobe.Cells.Item(93, 3).Value = 2771
$wsKobe.Cells.Item(93, 4).Value = 0
$wsKobe.Cells.Item(93, 5).Value = 281
$wsKobe.Cells.Item(93, 6).Value = 59
$wsKobe.Cells.Item(93, 7).Value = 50
$wsKobe.Cells.Item(93, 8).Value = 9
$wsKobe.Cells.Item(93, 9).Value = 11
$wsKobe.Cells.Item(93, 10).Value = 193

$wsKobe.Cells.Item(93, 1).Style = $wsKobe.Cells.Item(92, 1).Style
$wsKobe.Cells.Item(93, 2).Style = $wsKobe.Cells.Item(92, 2).Style
$wsKobe.Cells.Item(93, 3).Style = $wsKobe.Cells.Item(92, 3).Style
$wsKobe.Cells.Item(93, 4).Style = $wsKobe.Cells.Item(92, 4).Style
$wsKobe.Cells.Item(93, 5).Style = $wsKobe.Cells.Item(92, 5).Style
$wsKobe.Cells.Item(93, 6).Style = $wsKobe.Cells.Item(92, 6).Style
$wsKobe.Cells.Item(93, 7).Style = $wsKobe.Cells.Item(92, 7).Style
$wsKobe.Cells.Item(93, 8).Style = $wsKobe.Cells.Item(92, 8).Style
$wsKobe.Cells.Item(93, 9).Style = $wsKobe.Cells.Item(92, 9).Style
$wsKobe.Cells.Item(93, 10).Style = $wsKobe.Cells.Item(92, 10).Style

$wsKobe.Application.ActiveWindow.FreezePanes = $false
$wsKobe.Range("B85").Select()
$wsKobe.Application.ActiveWindow.FreezePanes = $true
$wsKobe.Range("A93").Select()

# ---------------------------------------------------------------------------
# Sheet "other" -> new data row 68 (was the "footnote" row), footnote pushes
# to row 69. This sheet becomes the "tabSelected" sheet.
# ---------------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")
$wsOther.Rows.Item(68).Insert()

$wsOther.Cells.Item(68, 1).Value = 43966
$wsOther.Cells.Item(68, 1).NumberFormat = "m/d/yyyy"
$wsOther.Cells.Item(68, 2).Value = 0
$wsOther.Cells.Item(68, 3).Value = 14
$wsOther.Cells.Item(68, 4).Value = 5
$wsOther.Cells.Item(68, 5).Value = 4
$wsOther.Cells.Item(68, 6).Value = 1
$wsOther.Cells.Item(68, 7).Value = 0
$wsOther.Cells.Item(68, 8).Value = 9

$wsOther.Cells.Item(68, 1).Style = $wsOther.Cells.Item(67, 1).Style
$wsOther.Cells.Item(68, 2).Style = $wsOther.Cells.Item(67, 2).Style
$wsOther.Cells.Item(68, 3).Style = $wsOther.Cells.Item(67, 3).Style
$wsOther.Cells.Item(68, 4).Style = $wsOther.Cells.Item(67, 4).Style
$wsOther.Cells.Item(68, 5).Style = $wsOther.Cells.Item(67, 5).Style
$wsOther.Cells.Item(68, 6).Style = $wsOther.Cells.Item(67, 6).Style
$wsOther.Cells.Item(68, 7).Style = $wsOther.Cells.Item(67, 7).Style
$wsOther.Cells.Item(68, 8).Style = $wsOther.Cells.Item(67, 8).Style

$wsOther.Range("A68").Select()

$wsOther.Select()

$wb.Windows.Item(1).ActiveTab = $wsOther
